$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 169117.17
$ws.Range("I6").Value = 202760.6
$ws.Range("K6").Value = 608281.8
$ws.Range("M6").Value = -608169.8
$ws.Range("H17").Value = 699.25714
$ws.Range("J17").Value = 711.08826
$ws.Range("L17").Value = 2133.26478
$ws.Range("N17").Value = -2469.26478
$ws.Range("H21").Value = 13960
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 13960
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 13960
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -14896
$ws.Range("H23").Value = 13960
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 13960
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 13960
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -14428
$ws.Range("H29").Value = 174.5
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H106").Value = 3075
$ws.Range("I106").Value = 3675
$ws.Range("J106").Value = 2475
$ws.Range("K106").Value = 3675
$ws.Range("L106").Value = 2475
$ws.Range("M106").Value = -3044
$ws.Range("N106").Value = -3737
$ws.Range("H137").Value = 1397.4117
$ws.Range("I137").Value = 1165.8276
$ws.Range("J137").Value = 2740.6
$ws.Range("K137").Value = 3497.4828
$ws.Range("L137").Value = 8221.799999999999
$ws.Range("M137").Value = -947.4828000000002
$ws.Range("N137").Value = -13321.8
$ws.Range("H138").Value = 5896.4585
$ws.Range("I138").Value = 1096
$ws.Range("J138").Value = 39499.668
$ws.Range("K138").Value = 3288
$ws.Range("L138").Value = 118499.004
$ws.Range("M138").Value = 1852
$ws.Range("N138").Value = -128779.004
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 734.7714
$ws.Range("I74").Value = 729.3929000000001
$ws.Range("J74").Value = 756.2857
$ws.Range("K74").Value = 729.3929000000001
$ws.Range("L74").Value = 756.2857
$ws.Range("M74").Value = 144.6070999999999
$ws.Range("N74").Value = -2504.2857
$ws.Range("H77").Value = 734.7714
$ws.Range("I77").Value = 729.3929000000001
$ws.Range("J77").Value = 756.2857
$ws.Range("K77").Value = 3646.9645
$ws.Range("L77").Value = 3781.4285
$ws.Range("M77").Value = 721.0355
$ws.Range("N77").Value = -12517.4285
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1378.4286
$ws.Range("I80").Value = 630.25
$ws.Range("J80").Value = 1838.8462
$ws.Range("K80").Value = 630.25
$ws.Range("L80").Value = 1838.8462
$ws.Range("M80").Value = 367.75
$ws.Range("N80").Value = -3834.8462
$ws.Range("H83").Value = 1378.4286
$ws.Range("I83").Value = 630.25
$ws.Range("J83").Value = 1838.8462
$ws.Range("K83").Value = 3151.25
$ws.Range("L83").Value = 9194.231
$ws.Range("M83").Value = 1840.75
$ws.Range("N83").Value = -19178.231
$ws.Range("H107").Value = 45493844
$ws.Range("I107").Value = 55603424
$ws.Range("K107").Value = 55603424
$ws.Range("M107").Value = -55601504
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1429714.2
$ws.Range("I4").Value = 3334000
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 3334000
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -3333888
$ws.Range("N4").Value = -1724
$ws.Range("H31").Value = 32897
$ws.Range("I31").Value = 1097.8667
$ws.Range("J31").Value = 48283.676
$ws.Range("K31").Value = 1097.8667
$ws.Range("L31").Value = 48283.676
$ws.Range("M31").Value = -802.8667
$ws.Range("N31").Value = -48873.676
$ws.Range("H34").Value = 32897
$ws.Range("I34").Value = 1097.8667
$ws.Range("J34").Value = 48283.676
$ws.Range("K34").Value = 1097.8667
$ws.Range("L34").Value = 48283.676
$ws.Range("M34").Value = -895.8667
$ws.Range("N34").Value = -48687.676
$ws.Range("H132").Value = 27780536
$ws.Range("I132").Value = 23258352
$ws.Range("K132").Value = 69775056
$ws.Range("M132").Value = -69772526
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 650.2069
$ws.Range("I113").Value = 616.25
$ws.Range("J113").Value = 674.17645
$ws.Range("K113").Value = 1848.75
$ws.Range("L113").Value = 2022.52935
$ws.Range("M113").Value = 321.25
$ws.Range("N113").Value = -6362.529350000001
$ws.Range("H124").Value = 4980
$ws.Range("J124").Value = 4980
$ws.Range("L124").Value = 14940
$ws.Range("N124").Value = -24760
$ws.Range("H125").Value = 781.6667
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 1820
$ws.Range("I126").Value = 1820
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5460
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -520
$ws.Range("N126").ClearContents()
$ws.Range("H131").Value = 8340.605
$ws.Range("I131").Value = 861.63635
$ws.Range("J131").Value = 9606.277
$ws.Range("K131").Value = 2584.90905
$ws.Range("L131").Value = 28818.831
$ws.Range("M131").Value = 2455.09095
$ws.Range("N131").Value = -38898.831
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 160
$ws.Range("I4").Value = 160
$ws.Range("K4").Value = 160
$ws.Range("M4").Value = -48
$ws.Range("H70").Value = 82047
$ws.Range("I70").Value = 147130.86
$ws.Range("J70").Value = 6115.8335
$ws.Range("K70").Value = 147130.86
$ws.Range("L70").Value = 6115.8335
$ws.Range("M70").Value = -146860.86
$ws.Range("N70").Value = -6655.8335
$ws.Range("H73").Value = 82047
$ws.Range("I73").Value = 147130.86
$ws.Range("J73").Value = 6115.8335
$ws.Range("K73").Value = 147130.86
$ws.Range("L73").Value = 6115.8335
$ws.Range("M73").Value = -146194.86
$ws.Range("N73").Value = -7987.8335
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 302600
$ws.Range("J2").Value = 6500
$ws.Range("L2").Value = 6500
$ws.Range("N2").Value = -6724
$ws.Range("H55").Value = 340.61765
$ws.Range("I55").Value = 210.6875
$ws.Range("J55").Value = 456.1111
$ws.Range("K55").Value = 210.6875
$ws.Range("L55").Value = 456.1111
$ws.Range("M55").Value = -37.6875
$ws.Range("N55").Value = -802.1111000000001
